# Add a new work-log entry row (row 19) and normalize the style of the
# preceding rows (14-18) so they match the rest of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Capture the "highlighted" style currently used by B14:C18 so it can be
#     reapplied later to the brand-new last row. -----------------------------
$ws.Range("B14:C14").Copy() | Out-Null
$ws.Range("B19:C19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Normalize styles of rows 14-18 to match rows 2-13 ----------------------
# Copy the style from B13/C13 (the "normal" style) down onto B14:C18 so the
# previously highlighted rows become regular rows again.
$ws.Range("B13:C13").Copy() | Out-Null
$ws.Range("B14:C18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Append the new row of data ---------------------------------------------
$ws.Range("A19").Value = 45706
$ws.Range("B19").Value = 2
$ws.Range("C19").Value = "python file to execute commands was created"

# Match the date-cell style used by the rest of column A.
$ws.Range("A18").Copy() | Out-Null
$ws.Range("A19").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# --- Adjust the view ---------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
